$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.387.39"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.876.87"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.59"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07960"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3150"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.98"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08137"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").Value = "1.887.92"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.34"
$ws.Range("E13").Value = "  +4.47%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7075"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.404"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008458"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "29.395.47"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.58"
$ws.Range("E19").Value = "  +5.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.42"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").Value = "2.138.34"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.675"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1588"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.065"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.82"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.507"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.421"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.300"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.218"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05327"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.947"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7569"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01892"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").Value = "1.268.66"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.763"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.395"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9068"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.29"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.71"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "2.033.10"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.812"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5204"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.519"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4347"
$ws.Range("E51").Value = "  -0.18%  "
